$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the row that will no longer be used (row 4)
$ws.Range("A4:D4").Clear()

# Header row
$ws.Range("A1").Value = "Plik"
$ws.Range("B1").Value = "Parametry"
$ws.Range("C1").Value = "Aktualne użycie pamięci"
$ws.Range("D1").Value = "Szczytowe użycie pamięci"

# Copy header style (already on A1:C1) to the new D1 cell
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 2
$ws.Range("A2").Value = "data/50_1.txt"
$ws.Range("B2").Value = "dobre"
$ws.Range("C2").Value = 0.3469924926757812
$ws.Range("D2").Value = 0.355438232421875

# Row 3
$ws.Range("A3").Value = "data/50_1.txt"
$ws.Range("B3").Value = "słabe"
$ws.Range("C3").Value = 0.06217193603515625
$ws.Range("D3").Value = 0.0703277587890625
